# Update "想去人数" (want-to-go count) figures that changed between the
# previous scrape and the latest gh-pages data refresh (commit 456a3b4).
#
# Sheet "展览" (Exhibition):
#   F2: 5268 -> 5274
#   F6:  798 ->  799
#   F7:  295 ->  298
#
# Sheet "演出" (Performance):
#   F2:   37 ->   39
#   F3:    9 ->   10
#
# Sheet "全部类型" (All types, combined roll-up of the two sheets above):
#   F2: 5268 -> 5274
#   F6:  798 ->  799
#   F7:   37 ->   39
#   F8:  295 ->  298
#   F10:   9 ->   10

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5274
$wsExhibition.Range("F6").Value = 799
$wsExhibition.Range("F7").Value = 298

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 39
$wsPerformance.Range("F3").Value = 10

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5274
$wsAll.Range("F6").Value = 799
$wsAll.Range("F7").Value = 39
$wsAll.Range("F8").Value = 298
$wsAll.Range("F10").Value = 10
